$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "Daniel Reichmann"

# Fill in row 2 data: date, from-time, to-time, category, comment
$ws.Range("A2").Value = 41666
$ws.Range("B2").Value = 0.79166666666666663
$ws.Range("C2").Value = 0.83333333333333337
$ws.Range("E2").Value = "UseCase"
$ws.Range("F2").Value = "Draft, needs aprovement"

# Update the selected cell to F2
$ws.Range("F2").Select()

$wb.Save()
